$d = $word.ActiveDocument

# Locate the empty paragraph that sits right after
# "Dat lai ten cho cac label , textbox" and right before
# "Trigger khi xoa khach hang ..." by searching for the start of the
# following paragraph's text, then counting paragraphs up to that
# point (Paragraph.Previous/.Next are unreliable in this host, so we
# navigate purely by paragraph index instead).
$hitRange = $d.Content
$found = $hitRange.Find.Execute("Trigger khi", $true, $false, $false, $false,
                                 $false, $true, 1, $false, "", 0)

$targetIndex = 142
if ($found) {
    $beforeRange = $d.Range(0, $hitRange.Start)
    $targetIndex = $beforeRange.Paragraphs.Count
}

$target = $d.Paragraphs.Item($targetIndex)

# Split the empty paragraph in two: the first stays empty (same
# paragraph mark formatting), the new second paragraph receives the
# note text, keeping identical run formatting (Times New Roman 14pt /
# sz 28 half-points, en-US language) already carried by the empty
# paragraph mark.
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.InsertAfter("Co the lam :")
